# Insert a new data row above row 35 (pushes old rows 35..118 down to 36..119)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(35).EntireRow.Insert()

# Populate the newly inserted row 35 with the new record.
$ws.Cells.Item(35, 1).Value = 10
$ws.Cells.Item(35, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(35, 3).Value = "La Araucanía"
$ws.Cells.Item(35, 4).Value = 45152
$ws.Cells.Item(35, 5).Value = 9
$ws.Cells.Item(35, 6).Value = "Fruta"
$ws.Cells.Item(35, 7).Value = 100108
$ws.Cells.Item(35, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(35, 9).Value = 100108007
$ws.Cells.Item(35, 10).Value = "Coco"
$ws.Cells.Item(35, 11).Value = "Sin especificar"
$ws.Cells.Item(35, 12).Value = "Primera"
$ws.Cells.Item(35, 13).Value = 35
$ws.Cells.Item(35, 14).Value = 36000
$ws.Cells.Item(35, 15).Value = 36000
$ws.Cells.Item(35, 16).Value = 36000
$ws.Cells.Item(35, 17).Value = "$/malla 20 unidades"
$ws.Cells.Item(35, 18).Value = "Perú"
$ws.Cells.Item(35, 19).Value = 1800
$ws.Cells.Item(35, 20).Value = 20
